$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column numbers (B,C,D,E,G,H,I,J,K,N,O) corresponding to data columns below
$cols = @(2, 3, 4, 5, 7, 8, 9, 10, 11, 14, 15)

$data = @(
    @(0.3714395591177038, 0.05579912950918242, 0.6523238794029282, 0.2662843189305661, 0.7169600431790144, 0.8196469962007171, 0.6284241979578198, 0.1387787624293466, 0.3991581964930049, 1.490711683824989, 3.075950561197459),
    @(0.3347798175139474, 0.04928632770985075, 0.6416451284709694, 0.2610608681683928, 0.7189266023718304, 0.8243504706904758, 0.6335247926551446, 0.1353281451824415, 0.35797459018994, 1.505711265271, 3.089697132661342),
    @(0.3123118185157239, 0.04528376544608648, 0.6354147157113914, 0.2579929372904743, 0.7205566758285613, 0.8275635345104391, 0.636962074204888, 0.1332862474594307, 0.332716901195738, 1.51540467414606, 3.099704559935248),
    @(0.3031668799846727, 0.04365183978664788, 0.6329580189775186, 0.2567778139994559, 0.7213271885480381, 0.8289547102847337, 0.6384396201625826, 0.1324734857647272, 0.3224321486333395, 1.519476445031003, 3.104176800790768),
    @(0.3016490487148076, 0.04338081071000488, 0.6325550590378555, 0.2565781643691167, 0.7214615478962116, 0.8291906582527275, 0.6386896058291107, 0.1323396951960092, 0.3207248722315228, 1.520159906197057, 3.104943220736715),
    @(0.3121884415497789, 0.04526176006334026, 0.6353812505923884, 0.2579764075838327, 0.720566637080978, 0.8275819649844323, 0.636981689833398, 0.1332752079887953, 0.332578164377594, 1.515459094975663, 3.099763278174294),
    @(0.3587910627427107, 0.05355432149008266, 0.6485742055535297, 0.264454396315486, 0.7175503882761518, 0.8212013214370586, 0.630119481093697, 0.137573059838914, 0.3849523824891037, 1.495783181383937, 3.0803652111383),
    @(0.450484579154363, 0.06978432140076052, 0.6770296115565202, 0.278261317091534, 0.71499038805716, 0.8112654664741399, 0.6190865270100652, 0.1466103070310893, 0.4878678449472886, 1.461034471380422, 3.054756504771802),
    @(0.5180164096261421, 0.08168694661624443, 0.6995065615402893, 0.289077331737893, 0.7151582708398792, 0.805532564705743, 0.6124582193425141, 0.1536220140453253, 0.5635855620276402, 1.437840124569208, 3.043519913930254),
    @(0.548769878548967, 0.0870966722000901, 0.7100722016656675, 0.2941436830106454, 0.7156803325401171, 0.8032640085847902, 0.6097636996144367, 0.1568928019712814, 0.5980502058388595, 1.427794969719839, 3.040053862600587),
    @(0.5604196376613402, 0.08914443758413881, 0.7141220026713881, 0.2960831483497941, 0.7159421613484511, 0.8024537024080018, 0.6087894791051589, 0.1581430259780063, 0.6111034549625742, 1.424063890614242, 3.038977958605642),
    @(0.5579104820284613, 0.08870345071881047, 0.713247637032822, 0.2956645190736893, 0.7158829186719089, 0.8026260491727868, 0.6089972427893393, 0.1578732501349975, 0.6082921155554004, 1.424864206804454, 3.039199150466175),
    @(0.5497282330003657, 0.087265159213473, 0.7104044037855601, 0.2943028245049888, 0.7157005877267011, 0.8031963675583995, 0.6096826251941572, 0.1569954253246095, 0.5991240635675297, 1.427486552329592, 3.039960605288911),
    @(0.5447168824831294, 0.08638406000886789, 0.70866919411867, 0.2934714735378705, 0.715597258419308, 0.8035520507087597, 0.6101084503734526, 0.156459248590906, 0.5935086381057033, 1.429102296744482, 3.040457832210819),
    @(0.5160072051700126, 0.08133330374269576, 0.698822913857839, 0.288749168411826, 0.7151331252274957, 0.8056876445294421, 0.6126407677411265, 0.1534098912851931, 0.5613335664415331, 1.438506785054376, 3.043779539890892),
    @(0.4984027170449963, 0.07823352900894065, 0.692869690364887, 0.2858895624444528, 0.7149625792709031, 0.8070846442131625, 0.6142764198583208, 0.1515599750003958, 0.5415999118916943, 1.444405771726803, 3.046238752101658),
    @(0.4882802163014617, 0.07645017021776823, 0.6894776386308763, 0.2842585471982915, 0.7149064335975481, 0.8079201077655682, 0.6152473914824945, 0.1505035892847673, 0.5302515846625226, 1.447846344512625, 3.047808111660174),
    @(0.4848534695498472, 0.07584628080194022, 0.6883346642577237, 0.2837086776178168, 0.7148946268261085, 0.8082084698326923, 0.6155813301194009, 0.1501472274840694, 0.526409591470383, 1.449019443234334, 3.048366071948436),
    @(0.5002764250364464, 0.07856355254331504, 0.6935001017899936, 0.286192549277267, 0.7149763924017662, 0.8069326252834514, 0.6140991774667093, 0.1517561114580275, 0.5437003963169218, 1.443772884609221, 3.045960934745182),
    @(0.5521314520640317, 0.08768764214408975, 0.7112382067153931, 0.2947022191620334, 0.7157524017685972, 0.8030275286795074, 0.6094800595405481, 0.1572529477857074, 0.6018168873872014, 1.426714329770549, 3.03973052598937),
    @(0.5860452385384747, 0.09364617476899184, 0.7231155743461386, 0.3003858536407265, 0.7166334195592299, 0.800759431080138, 0.606730119723526, 0.1609133304536243, 0.6398120946615506, 1.415989945970985, 3.037037758426493),
    @(0.5679428721049931, 0.09046644145433902, 0.7167504281218271, 0.2973412411144949, 0.7161289814915079, 0.8019439792658858, 0.6081732031295282, 0.1589535102611705, 0.6195323979098646, 1.421674910890708, 3.038348749130591),
    @(0.4994293261090945, 0.07841435281328302, 0.6932149976571793, 0.2860555284009081, 0.7149700169619422, 0.8070012523935333, 0.6141792133936868, 0.1516674158755364, 0.5427507771824764, 1.444058859724025, 3.046086051515829),
    @(0.4256485661184115, 0.06539730009923517, 0.6690554964424962, 0.274408103330714, 0.7153234121859384, 0.8136779312411733, 0.6218117499696305, 0.1441002053231628, 0.4600063177901177, 1.470024581876679, 3.060353537645767)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rowdata = $data[$i]
    for ($j = 0; $j -lt $cols.Length; $j++) {
        $ws.Cells.Item($row, $cols[$j]).Value = $rowdata[$j]
    }
}

Write-Output "Updated pl_mw sheet values for 380 kV case"